# 17.06.19 Today Sales Updated
#
# Updates the "Raju Ahamed" sheet with the new day's cash count:
#  - Both date headers (Today table + This Month table) move from
#    16.06.19 to 17.06.19.
#  - The previous "Today" block is marked with a "\" so it is visually
#    struck through / closed out once the new figures are entered.
#  - The Today denomination quantities (and their dependent Amount /
#    Grand Total formulas) are refreshed with today's counts.
#  - Because this is the first day recorded for the month, the "This
#    Month" table mirrors the same quantities as "Today".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raju Ahamed")

# Update the date shown above both tables on the sheet.
$ws.Range("A4").Value  = "Date: 17.06.19"
$ws.Range("A31").Value = "Date: 17.06.19"

# Flag the old "Today" section as closed out.
$ws.Range("A32").Value = "\"

# Today table (rows 6-11): Qty. column per denomination (1000,500,100,50,20,10)
$ws.Range("E6").Value  = 17
$ws.Range("E7").Value  = 60
$ws.Range("E8").Value  = 392
$ws.Range("E9").Value  = 2
$ws.Range("E10").Value = 135
$ws.Range("E11").Value = ""

# This Month table (rows 33-38): same denomination quantities as Today
$ws.Range("E33").Value = 17
$ws.Range("E34").Value = 60
$ws.Range("E35").Value = 392
$ws.Range("E36").Value = 2
$ws.Range("E37").Value = 135
$ws.Range("E38").Value = ""
